$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at S (shifts old S->T, T->U, U->V, V->W for cell data/styles).
$ws.Columns("S").Insert()

# New header for the inserted "Agreement Unit Type" column.
$ws.Range("S1").Value = "Agreement Unit Type"

# New data for the inserted column.
$ws.Range("S2").Value = "A1"
$ws.Range("S3").Value = "A2"
$ws.Range("S4").Value = "A3"
$ws.Range("S5").Value = "A4"

# Comments on the header row are anchored to the original cell references and are not
# shifted automatically by the column insert, so fix them up to match the new layout.
# R1 (Agreement Committed Amount) and S1 (now Agreement Unit Type) keep their existing
# comments as-is; T1 already carries the correct "SEBI reporting" text. U1 and V1 need to
# be updated to hold what used to live in T1/U1 respectively.
$null = $ws.Range("U1").Comment.Text("Author:`n-Mandatory`n-As per SEBI reporting requirements")
$ws.Range("V1").AddComment("Author:`n-Incase custom fields are already created, please download the data from the platform to use existing headers to avoid errors while generating documents`n-Can add any custom fields as per your requirement") | Out-Null
